$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Input")
$ws2 = $wb.Worksheets.Item("JSON")

# Simulate a refresh of the "ContentCache" Power Query: a new child node
# ("A.A", id 1005, parent 1001) was added to the Input table, so the table
# grows by one row and the cached JSON blob on the JSON sheet is rebuilt
# (picking up fresh GUIDs for every node along the way).

$tbl = $ws1.ListObjects.Item(1)
$tbl.ListRows.Add() | Out-Null

$ws1.Range("B8").Value = 1005
$ws1.Range("C8").Value = 1001
$ws1.Range("D8").Value = 0
$ws1.Range("E8").Value = 1057
$ws1.Range("F8").Value = "A.A"
$ws1.Range("G8").Value = 1056
$ws1.Range("H8").Value = "Here's a nice heading"
$ws1.Range("I8").Value = "<p>Here's a summary</p>"

$ws2.Range("A2").Value = '{"1000":{"Node":{"Id":1000,"ParentContentId":-1,"SortOrder":0,"Uid":"778ab36e-52af-4498-b1dc-6668ec00be12"},"ContentTypeId":1057,"PublishedData":{"Name":"Site","UrlSegment":"site","TemplateId":1056,"Published":true,"Properties":{"listHeading":[{"v":null}],"listSummary":[{"v":null}]}}},"1001":{"Node":{"Id":1001,"ParentContentId":1000,"SortOrder":0,"Uid":"7ce2bf02-c1d2-4bd0-b796-39403f4a0c70"},"ContentTypeId":1057,"PublishedData":{"Name":"A","UrlSegment":"a","TemplateId":1056,"Published":true,"Properties":{"listHeading":[{"v":"A"}],"listSummary":[{"v":null}]}}},"1002":{"Node":{"Id":1002,"ParentContentId":1000,"SortOrder":1,"Uid":"b00581e5-0ec9-4c58-94d5-176d1a098765"},"ContentTypeId":1057,"PublishedData":{"Name":"B","UrlSegment":"b","TemplateId":1056,"Published":true,"Properties":{"listHeading":[{"v":"B"}],"listSummary":[{"v":null}]}}},"1003":{"Node":{"Id":1003,"ParentContentId":1000,"SortOrder":2,"Uid":"2f854b70-5bec-4b07-8481-a5e7ac314d6c"},"ContentTypeId":1057,"PublishedData":{"Name":"C","UrlSegment":"c","TemplateId":1056,"Published":true,"Properties":{"listHeading":[{"v":"C"}],"listSummary":[{"v":null}]}}},"1004":{"Node":{"Id":1004,"ParentContentId":1000,"SortOrder":3,"Uid":"4243b302-1ac3-45fd-ac92-cac3b2f88a98"},"ContentTypeId":1057,"PublishedData":{"Name":"D","UrlSegment":"d","TemplateId":1056,"Published":true,"Properties":{"listHeading":[{"v":"D"}],"listSummary":[{"v":null}]}}},"1005":{"Node":{"Id":1005,"ParentContentId":1001,"SortOrder":0,"Uid":"1a1e032d-8c8a-4a5f-b834-6005762314ab"},"ContentTypeId":1057,"PublishedData":{"Name":"A.A","UrlSegment":"aa","TemplateId":1056,"Published":true,"Properties":{"listHeading":[{"v":"Here''s a nice heading"}],"listSummary":[{"v":"<p>Here''s a summary<\/p>"}]}}}}'

$ws1.Activate()
$ws1.Range("I7").Select()
